$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert three new date columns before column B. This shifts the
#    existing B,C,D,E columns (and their values/styles) to E,F,G,H.
$ws.Range("B1:D1").EntireColumn.Insert()

# 2. New header values for the inserted date columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# 3. The newly inserted columns are blank - fill them with the same
#    default "UN" marker used across the rest of the grid (rows 2-27).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# 4. Re-apply the explicit column widths (8 characters) across all the
#    date columns C:H, keeping the same look as before the insert.
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 7.2
$ws.Range("H1").EntireColumn.ColumnWidth = 7.2

# 5. Append the two new analyst/benchmark rows at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
